$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 35 (shifts rows 35+ down by one, and Excel
# automatically re-points the SUM() ranges that referenced row 40 as the
# last data row, e.g. SUM(C6:C40) -> SUM(C6:C41)).
$ws.Rows.Item(35).EntireRow.Insert()

# Fill in the data for the newly inserted row (new task logged).
$ws.Range("A35").Value = "Integration database and network"
$ws.Range("X35").Value = 2
$ws.Range("Y35").Value = 4

# New hours logged against existing tasks for Wednesday 19th March (column Y).
$ws.Range("Y26").Value = 6
$ws.Range("Y31").Value = 14
$ws.Range("Y32").Value = 6
$ws.Range("Y33").Value = 6

# Extend the weekly total row's shared formula/number formatting out to the
# two newly-used days (X = 18th, Y = 19th) as well as the still-blank Z/AA.
$ws.Range("X3:AA3").NumberFormat = "0.00"
$ws.Range("X3").Formula = "=SUM(X6:X41)"
$ws.Range("Y3").Formula = "=SUM(Y6:Y41)"
$ws.Range("Z3").Formula = "=SUM(Z6:Z41)"
$ws.Range("AA3").Formula = "=SUM(AA6:AA41)"

# Move the active selection to reflect where the user was last working.
$ws.Range("A27").Select()
